# Updated cryptos list on Mon Sep  4 13:35:58 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.914.72"
$ws.Range("E2").Value = "  -0.19%  "
$ws.Range("D3").Value = "1.636.06"
$ws.Range("E3").Value = "  -0.36%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  -0.33%  "
$ws.Range("D5").Value = "'216.74"
$ws.Range("E5").Value = "  +0.74%  "
$ws.Range("D6").Value = "'0.5063"
$ws.Range("E6").Value = "  -0.23%  "
$ws.Range("D7").Value = "'1.002"
$ws.Range("E7").Value = "  -0.18%  "
$ws.Range("D8").Value = "'0.2581"
$ws.Range("E8").Value = "  +0.81%  "
$ws.Range("D9").Value = "'0.06366"
$ws.Range("E9").Value = "  -0.08%  "
$ws.Range("D10").Value = "'19.66"
$ws.Range("E10").Value = "  +0.66%  "
$ws.Range("D11").Value = "'0.07756"
$ws.Range("E11").Value = "  -0.17%  "
$ws.Range("E12").Value = "  -0.40%  "
$ws.Range("D13").Value = "1.629.47"
$ws.Range("E13").Value = "  -0.69%  "
$ws.Range("D14").Value = "'0.5494"
$ws.Range("E14").Value = "  +0.70%  "
$ws.Range("D15").Value = "0.0₅7740"
$ws.Range("E15").Value = "  -1.10%  "
$ws.Range("D16").Value = "'64.17"
$ws.Range("E16").Value = "  -0.18%  "
$ws.Range("D17").Value = "25.896.09"
$ws.Range("E17").Value = "  -0.39%  "
$ws.Range("E18").Value = "  -0.25%  "
$ws.Range("D19").Value = "'4.446"
$ws.Range("E19").Value = "  +0.12%  "
$ws.Range("D20").Value = "'195.14"
$ws.Range("E20").Value = "  -1.16%  "
$ws.Range("D21").Value = "'9.907"
$ws.Range("E21").Value = "  -0.58%  "
$ws.Range("D22").Value = "'6.084"
$ws.Range("E22").Value = "  +0.55%  "
$ws.Range("E23").Value = "  -0.26%  "
$ws.Range("E24").Value = "  +0.63%  "
$ws.Range("D25").Value = "'142.51"
$ws.Range("E25").Value = "  +0.74%  "
$ws.Range("D26").Value = "'0.1240"
$ws.Range("E26").Value = "  +5.78%  "
$ws.Range("D27").Value = "'6.831"
$ws.Range("E27").Value = "  -0.75%  "
$ws.Range("D28").Value = "'15.62"
$ws.Range("E28").Value = "  -0.67%  "
$ws.Range("D29").Value = "'1.246"
$ws.Range("E29").Value = "  +0.66%  "
$ws.Range("D30").Value = "'0.04862"
$ws.Range("E30").Value = "  -2.69%  "
$ws.Range("D31").Value = "'3.240"
$ws.Range("E31").Value = "  -0.64%  "
$ws.Range("D32").Value = "'3.195"
$ws.Range("E32").Value = "  +0.20%  "
$ws.Range("D33").Value = "'1.545"
$ws.Range("E33").Value = "  +0.22%  "
$ws.Range("D34").Value = "'2.369"
$ws.Range("E34").Value = "  +0.39%  "
$ws.Range("D35").Value = "'0.9052"
$ws.Range("E35").Value = "  +1.20%  "
$ws.Range("D36").Value = "'2.573"
$ws.Range("E36").Value = "  -0.49%  "

# Rows 37 and 38 swapped coin (Maker <-> ImmutableX) with new price/volume values
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").Value = "'0.5498"
$ws.Range("E37").Value = "  +0.90%  "

$ws.Range("B38").Value = "Maker"
$ws.Range("C38").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D38").Value = "1.124.59"
$ws.Range("E38").Value = "  -0.56%  "

$ws.Range("D39").Value = "'0.01557"
$ws.Range("E39").Value = "  +0.14%  "
$ws.Range("D40").Value = "'1.002"
$ws.Range("E40").Value = "  -0.24%  "
$ws.Range("D41").Value = "'5.579"
$ws.Range("E41").Value = "  -0.34%  "
$ws.Range("D42").Value = "'0.8047"
$ws.Range("E42").Value = "  -1.72%  "
$ws.Range("D43").Value = "'97.89"
$ws.Range("E43").Value = "  -1.93%  "
$ws.Range("D44").Value = "0.0₈121"
$ws.Range("E44").Value = "  -5.60%  "
$ws.Range("D45").Value = "1.770.10"
$ws.Range("E45").Value = "  -0.40%  "
$ws.Range("D46").Value = "'0.4458"
$ws.Range("E46").Value = "  -1.82%  "
$ws.Range("D47").Value = "'0.9993"
$ws.Range("E47").Value = "  -0.40%  "
$ws.Range("D48").Value = "'54.85"
$ws.Range("E48").Value = "  +0.04%  "
$ws.Range("D49").Value = "'0.05153"
$ws.Range("E49").Value = "  +1.53%  "
$ws.Range("D50").Value = "'7.487"
$ws.Range("E50").Value = "  +1.75%  "
$ws.Range("E51").Value = "  +0.01%  "
